$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values for rows 2-17 (natmi re-run per commit: "Natmi following Dr Hou advice")
# Columns: E,K = expressing-cell counts; G/H/I/J, M/N/O/P = ligand/receptor stats; Q/R/S/T = edge stats
$data = @{}
$data["E"] = @(3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3)
$data["G"] = @(29.82164133333333, 29.82164133333333, 29.82164133333333, 29.82164133333333, 1160.126729666667, 1160.126729666667, 1160.126729666667, 1160.126729666667, 0.2607933333333334, 0.2607933333333334, 0.2607933333333334, 0.2607933333333334, 86.94000199999999, 86.94000199999999, 86.94000199999999, 86.94000199999999)
$data["H"] = @(89.464924, 89.464924, 89.464924, 89.464924, 3480.380189, 3480.380189, 3480.380189, 3480.380189, 0.7823800000000001, 0.7823800000000001, 0.7823800000000001, 0.7823800000000001, 260.820006, 260.820006, 260.820006, 260.820006)
$data["I"] = @(0.02335016309719764, 0.02335016309719764, 0.02335016309719764, 0.02335016309719764, 0.9083721465342723, 0.9083721465342723, 0.9083721465342723, 0.9083721465342723, 0.0002041995878070102, 0.0002041995878070102, 0.0002041995878070102, 0.0002041995878070102, 0.06807349078072281, 0.06807349078072281, 0.06807349078072281, 0.06807349078072281)
$data["J"] = @(0.02335016309719765, 0.02335016309719765, 0.02335016309719765, 0.02335016309719765, 0.9083721465342726, 0.9083721465342726, 0.9083721465342726, 0.9083721465342726, 0.0002041995878070102, 0.0002041995878070102, 0.0002041995878070102, 0.0002041995878070102, 0.06807349078072282, 0.06807349078072282, 0.06807349078072282, 0.06807349078072282)
$data["K"] = @(3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3, 3)
$data["M"] = @(8.131233999999999, 243.3763986666667, 103.9426383333333, 33.41874933333333, 8.131233999999999, 243.3763986666667, 103.9426383333333, 33.41874933333333, 8.131233999999999, 243.3763986666667, 103.9426383333333, 33.41874933333333, 8.131233999999999, 243.3763986666667, 103.9426383333333, 33.41874933333333)
$data["N"] = @(24.393702, 730.1291960000001, 311.827915, 100.256248, 24.393702, 730.1291960000001, 311.827915, 100.256248, 24.393702, 730.1291960000001, 311.827915, 100.256248, 24.393702, 730.1291960000001, 311.827915, 100.256248)
$data["O"] = @(0.02090995573015822, 0.625857000534647, 0.2672947262403034, 0.08593831749489127, 0.02090995573015822, 0.625857000534647, 0.2672947262403034, 0.08593831749489127, 0.02090995573015822, 0.625857000534647, 0.2672947262403034, 0.08593831749489127, 0.02090995573015822, 0.625857000534647, 0.2672947262403034, 0.08593831749489127)
$data["P"] = @(0.02090995573015823, 0.6258570005346471, 0.2672947262403035, 0.08593831749489128, 0.02090995573015823, 0.6258570005346471, 0.2672947262403035, 0.08593831749489128, 0.02090995573015823, 0.6258570005346471, 0.2672947262403035, 0.08593831749489128, 0.02090995573015823, 0.6258570005346471, 0.2672947262403035, 0.08593831749489128)
$data["Q"] = @(242.4867439454053, 7257.883670035679, 3099.740079617051, 996.601956427239, 9433.261908574406, 282347.4654632109, 120586.6330825751, 38769.98437363009, 2.120571618973333, 63.47094226294224, 27.10754712641112, 8.715387034471112, 706.9295002224678, 21159.1445868328, 9036.773184585276, 2905.426133877498)
$data["R"] = @(2182.380695508647, 65320.95303032111, 27897.66071655346, 8969.417607845151, 84899.35717716966, 2541127.189168898, 1085279.697743176, 348929.8593626708, 19.08514457076, 571.2384803664801, 243.9679241377, 78.43848331024, 6362.365502002211, 190432.3012814952, 81330.95866126749, 26148.83520489748)
$data["S"] = @(0.0004882508766543768, 0.01461386303800692, 0.006241375452731878, 0.002006673729804464, 0.01899402137054043, 0.5685110669991585, 0.2428030842321951, 0.07806397393237818, 0.00000426980434116114, 0.0001277997415353067, 0.00005458147292125759, 0.00001754856900928477, 0.001423413678622248, 0.04260427075594612, 0.01819568508245512, 0.00585012126369931)
$data["T"] = @(0.0004882508766543771, 0.01461386303800692, 0.006241375452731882, 0.002006673729804465, 0.01899402137054044, 0.5685110669991588, 0.2428030842321952, 0.07806397393237823, 0.000004269804341161141, 0.0001277997415353067, 0.00005458147292125761, 0.00001754856900928477, 0.001423413678622248, 0.04260427075594614, 0.01819568508245513, 0.005850121263699312)

for ($row = 2; $row -le 17; $row++) {
    $idx = $row - 2
    foreach ($col in $data.Keys) {
        $ws.Range("$col$row").Value = $data[$col][$idx]
    }
}
